# Table 4.2 monthly update: October 2016 -> November 2016
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Title text and "Rolling 12 months" caption: October -> November
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Table 4.2. Receipts, Average Cost, and Quality of Fossil Fuels: Electric Utilities, 2006 - November 2016 (continued)"

# ------------------------------------------------------------------
# 2. Insert a new row 53 (November 2016 data) into the "Natural Gas"
#    block, pushing the annual-total / rolling-12-months rows down by
#    one. Copy the format of the preceding month row (52 = October)
#    into the freshly inserted row so the style indices line up with
#    the other month rows.
# ------------------------------------------------------------------
$ws.Range("A53:M53").EntireRow.Insert()
$ws.Range("A52:M52").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 7871
$ws.Range("C53").Value = 279
$ws.Range("D53").Value = 2.22
$ws.Range("E53").Value = 62.85
$ws.Range("F53").Value = 5.74
$ws.Range("G53").Value = 116.3
$ws.Range("H53").Value = 338187
$ws.Range("I53").Value = 326505
$ws.Range("J53").Value = 3.37
$ws.Range("K53").Value = 3.49
$ws.Range("L53").Value = 98
$ws.Range("M53").Value = 2.54

# ------------------------------------------------------------------
# 3. Update "Year to Date" annual totals (now rows 55-57) to reflect
#    data through November 2016.
# ------------------------------------------------------------------
# Year 2014 (row 55)
$ws.Range("B55").Value = 109499
$ws.Range("C55").Value = 3849
$ws.Range("D55").Value = 1.89
$ws.Range("E55").Value = 53.69
$ws.Range("F55").Value = 5.58
$ws.Range("G55").Value = 123.8
$ws.Range("H55").Value = 3573515
$ws.Range("I55").Value = 3478771
$ws.Range("J55").Value = 5.22
$ws.Range("K55").Value = 5.36
$ws.Range("L55").Value = 96.8
$ws.Range("M55").Value = 3.16

# Year 2015 (row 56)
$ws.Range("B56").Value = 107437
$ws.Range("C56").Value = 3772
$ws.Range("D56").Value = 1.8
$ws.Range("E56").Value = 51.37
$ws.Range("F56").Value = 5.19
$ws.Range("G56").Value = 130.2
$ws.Range("H56").Value = 4331629
$ws.Range("I56").Value = 4191468
$ws.Range("J56").Value = 3.57
$ws.Range("K56").Value = 3.69
$ws.Range("L56").Value = 96.2
$ws.Range("M56").Value = 2.69

# Year 2016 (row 57)
$ws.Range("B57").Value = 91688
$ws.Range("C57").Value = 3254
$ws.Range("D57").Value = 1.48
$ws.Range("E57").Value = 41.6
$ws.Range("F57").Value = 5.38
$ws.Range("G57").Value = 102.6
$ws.Range("H57").Value = 4702060
$ws.Range("I57").Value = 4547331
$ws.Range("J57").Value = 3.08
$ws.Range("K57").Value = 3.18
$ws.Range("L57").Value = 96.5
$ws.Range("M57").Value = 2.52

# ------------------------------------------------------------------
# 4. "Rolling 12 Months Ending in October" -> "...November" (row 58),
#    and update the rolling-12-months totals (now rows 59-60).
# ------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# Year 2015 rolling (row 59)
$ws.Range("B59").Value = 121731
$ws.Range("C59").Value = 4272
$ws.Range("D59").Value = 1.81
$ws.Range("E59").Value = 51.72
$ws.Range("F59").Value = 5.21
$ws.Range("G59").Value = 132.2
$ws.Range("H59").Value = 4634663
$ws.Range("I59").Value = 4485293
$ws.Range("J59").Value = 3.64
$ws.Range("K59").Value = 3.76
$ws.Range("L59").Value = 96.2
$ws.Range("M59").Value = 2.73

# Year 2016 rolling (row 60)
$ws.Range("B60").Value = 100180
$ws.Range("C60").Value = 3552
$ws.Range("D60").Value = 1.46
$ws.Range("E60").Value = 41.35
$ws.Range("F60").Value = 5.41
$ws.Range("G60").Value = 104.3
$ws.Range("H60").Value = 5088179
$ws.Range("I60").Value = 4920903
$ws.Range("J60").Value = 3.07
$ws.Range("K60").Value = 3.17
$ws.Range("L60").Value = 96.4
$ws.Range("M60").Value = 2.52

Write-Host "done"
